$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1131.3334
$ws.Range("I6").Value = 1169.2858
$ws.Range("J6").Value = 998.5
$ws.Range("K6").Value = 3507.8574
$ws.Range("L6").Value = 2995.5
$ws.Range("M6").Value = -3395.8574
$ws.Range("N6").Value = -3219.5
$ws.Range("H28").Value = 1003.5294
$ws.Range("I28").Value = 1003.5294
$ws.Range("K28").Value = 1003.5294
$ws.Range("M28").Value = -518.5294
$ws.Range("H87").Value = 19999.908
$ws.Range("J87").Value = 19999.908
$ws.Range("L87").Value = 19999.908
$ws.Range("N87").Value = -22495.908
$ws.Range("H90").Value = 19999.908
$ws.Range("J90").Value = 19999.908
$ws.Range("L90").Value = 59999.724
$ws.Range("N90").Value = -72479.724
$ws.Range("H98").Value = 1183.5
$ws.Range("I98").Value = 911.3333
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 911.3333
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 586.6667
$ws.Range("N98").Value = -4996
$ws.Range("H103").Value = 1374.75
$ws.Range("I103").Value = 749.5
$ws.Range("K103").Value = 2248.5
$ws.Range("M103").Value = -1662.5
$ws.Range("H107").Value = 497.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 497.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 497.5
$ws.Range("N107").Value = -4337.5
$ws.Range("M107").ClearContents()
$ws.Range("H122").Value = 1183.5
$ws.Range("I122").Value = 911.3333
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2733.9999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -283.9998999999998
$ws.Range("N122").Value = -10900
$ws.Range("H138").Value = 6176054
$ws.Range("I138").Value = 1213.8889
$ws.Range("J138").Value = 9263474
$ws.Range("K138").Value = 3641.6667
$ws.Range("L138").Value = 27790422
$ws.Range("M138").Value = 1498.3333
$ws.Range("N138").Value = -27800702

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17564.727
$ws.Range("I32").Value = 14892.897
$ws.Range("K32").Value = 14892.897
$ws.Range("M32").Value = -14605.897
$ws.Range("H97").Value = 905.931
$ws.Range("I97").Value = 968.1539
$ws.Range("J97").Value = 366.66666
$ws.Range("K97").Value = 968.1539
$ws.Range("L97").Value = 366.66666
$ws.Range("M97").Value = -472.1539
$ws.Range("N97").Value = -1358.66666
$ws.Range("H110").Value = 4662.222
$ws.Range("I110").Value = 4662.222
$ws.Range("K110").Value = 4662.222
$ws.Range("M110").Value = -2617.222
$ws.Range("H132").Value = 3206.2632
$ws.Range("I132").Value = 2939.9443
$ws.Range("K132").Value = 8819.832900000001
$ws.Range("M132").Value = -6289.832900000001
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 172.42105
$ws.Range("I80").Value = 67
$ws.Range("J80").Value = 200.53334
$ws.Range("K80").Value = 67
$ws.Range("L80").Value = 200.53334
$ws.Range("M80").Value = 931
$ws.Range("N80").Value = -2196.53334
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 172.42105
$ws.Range("I83").Value = 67
$ws.Range("J83").Value = 200.53334
$ws.Range("K83").Value = 335
$ws.Range("L83").Value = 1002.6667
$ws.Range("M83").Value = 4657
$ws.Range("N83").Value = -10986.6667
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 1544.175
$ws.Range("I107").Value = 1517.1143
$ws.Range("K107").Value = 1517.1143
$ws.Range("M107").Value = 402.8857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 977.5
$ws.Range("I16").Value = 553.3333
$ws.Range("K16").Value = 553.3333
$ws.Range("M16").Value = -266.3333
$ws.Range("H22").Value = 294287.06
$ws.Range("I22").Value = 206
$ws.Range("J22").Value = 500143.8
$ws.Range("K22").Value = 206
$ws.Range("L22").Value = 500143.8
$ws.Range("M22").Value = 144
$ws.Range("N22").Value = -500843.8
$ws.Range("H31").Value = 3469.5
$ws.Range("I31").Value = 2463.5
$ws.Range("J31").Value = 4307.8335
$ws.Range("K31").Value = 2463.5
$ws.Range("L31").Value = 4307.8335
$ws.Range("M31").Value = -2168.5
$ws.Range("N31").Value = -4897.8335
$ws.Range("H34").Value = 3469.5
$ws.Range("I34").Value = 2463.5
$ws.Range("J34").Value = 4307.8335
$ws.Range("K34").Value = 2463.5
$ws.Range("L34").Value = 4307.8335
$ws.Range("M34").Value = -2261.5
$ws.Range("N34").Value = -4711.8335
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H52").Value = 67216.75
$ws.Range("J52").Value = 67216.75
$ws.Range("L52").Value = 67216.75
$ws.Range("N52").Value = -67804.75
$ws.Range("H105").Value = 955.1818
$ws.Range("I105").Value = 745
$ws.Range("K105").Value = 745
$ws.Range("M105").Value = 1002
$ws.Range("H113").Value = 977.5
$ws.Range("I113").Value = 553.3333
$ws.Range("K113").Value = 553.3333
$ws.Range("M113").Value = 1616.6667
$ws.Range("H122").Value = 343271.53
$ws.Range("I122").Value = 928853.2
$ws.Range("K122").Value = 2786559.6
$ws.Range("M122").Value = -2784109.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 88074856
$ws.Range("I4").Value = 80100456
$ws.Range("J4").Value = 111112000
$ws.Range("K4").Value = 240301368
$ws.Range("L4").Value = 333336000
$ws.Range("M4").Value = -240301256
$ws.Range("N4").Value = -333336224
$ws.Range("H12").Value = 163.35715
$ws.Range("J12").Value = 178.60869
$ws.Range("L12").Value = 535.82607
$ws.Range("N12").Value = -881.82607
$ws.Range("H23").Value = 212.33333
$ws.Range("J23").Value = 247.4
$ws.Range("L23").Value = 742.2
$ws.Range("N23").Value = -1212.2
$ws.Range("H33").Value = 197.4
$ws.Range("J33").Value = 576
$ws.Range("L33").Value = 3456
$ws.Range("N33").Value = -4022
$ws.Range("H38").Value = 1490.2142
$ws.Range("J38").Value = 2074.8
$ws.Range("L38").Value = 6224.400000000001
$ws.Range("N38").Value = -6918.400000000001
$ws.Range("H74").Value = 6664.5835
$ws.Range("J74").Value = 6664.5835
$ws.Range("L74").Value = 19993.7505
$ws.Range("N74").Value = -22115.7505
$ws.Range("H77").Value = 6664.5835
$ws.Range("J77").Value = 6664.5835
$ws.Range("L77").Value = 59981.2515
$ws.Range("N77").Value = -70589.2515
$ws.Range("H113").Value = 1416.0454
$ws.Range("I113").Value = 1326.4445
$ws.Range("J113").Value = 1478.0769
$ws.Range("K113").Value = 3979.3335
$ws.Range("L113").Value = 4434.2307
$ws.Range("M113").Value = -1809.3335
$ws.Range("N113").Value = -8774.2307
$ws.Range("H114").Value = 1745.25
$ws.Range("I114").Value = 279
$ws.Range("K114").Value = 837
$ws.Range("M114").Value = 2417
$ws.Range("H122").Value = 1890.7222
$ws.Range("J122").Value = 2052.7273
$ws.Range("L122").Value = 18474.5457
$ws.Range("N122").Value = -23374.5457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 28021.568
$ws.Range("I80").Value = 40876.656
$ws.Range("K80").Value = 40876.656
$ws.Range("M80").Value = -39878.656
$ws.Range("H83").Value = 28021.568
$ws.Range("I83").Value = 40876.656
$ws.Range("K83").Value = 204383.28
$ws.Range("M83").Value = -199391.28
$ws.Range("H132").Value = 3598.5386
$ws.Range("I132").Value = 2979.4546
$ws.Range("J132").Value = 7003.5
$ws.Range("K132").Value = 8938.363799999999
$ws.Range("L132").Value = 21010.5
$ws.Range("M132").Value = -6408.363799999999
$ws.Range("N132").Value = -26070.5
$ws.Range("H136").Value = 16977.562
$ws.Range("J136").Value = 16977.562
$ws.Range("L136").Value = 50932.686
$ws.Range("N136").Value = -56032.686

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 25749.75
$ws.Range("J6").Value = 25749.75
$ws.Range("L6").Value = 25749.75
$ws.Range("N6").Value = -25973.75
$ws.Range("H68").Value = 6684.1304
$ws.Range("I68").Value = 4556.4287
$ws.Range("J68").Value = 7615
$ws.Range("K68").Value = 4556.4287
$ws.Range("L68").Value = 7615
$ws.Range("M68").Value = -3807.4287
$ws.Range("N68").Value = -9113
$ws.Range("H71").Value = 6684.1304
$ws.Range("I71").Value = 4556.4287
$ws.Range("J71").Value = 7615
$ws.Range("K71").Value = 22782.1435
$ws.Range("L71").Value = 38075
$ws.Range("M71").Value = -19038.1435
$ws.Range("N71").Value = -45563
$ws.Range("H102").Value = 47599.2
$ws.Range("J102").Value = 47599.2
$ws.Range("L102").Value = 47599.2
$ws.Range("N102").Value = -54089.2
$ws.Range("H104").Value = 21146.428
$ws.Range("J104").Value = 21146.428
$ws.Range("L104").Value = 21146.428
$ws.Range("N104").Value = -28134.428
$ws.Range("H105").Value = 49807.5
$ws.Range("J105").Value = 49807.5
$ws.Range("L105").Value = 49807.5
$ws.Range("N105").Value = -56795.5
$ws.Range("H122").Value = 3944.4878
$ws.Range("I122").Value = 3165.2058
$ws.Range("K122").Value = 9495.617400000001
$ws.Range("M122").Value = -7045.617400000001
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 4127.2
$ws.Range("I136").Value = 2232.8696
$ws.Range("K136").Value = 6698.6088
$ws.Range("M136").Value = -4148.6088
